$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 201297
$ws.Range("I70").Value = 1621.25
$ws.Range("J70").Value = 1000000
$ws.Range("K70").Value = 4863.75
$ws.Range("L70").Value = 3000000
$ws.Range("M70").Value = -4593.75
$ws.Range("N70").Value = -3000540

$ws.Range("H73").Value = 201297
$ws.Range("I73").Value = 1621.25
$ws.Range("J73").Value = 1000000
$ws.Range("K73").Value = 4863.75
$ws.Range("L73").Value = 3000000
$ws.Range("M73").Value = -3927.75
$ws.Range("N73").Value = -3001872

$ws.Range("H106").Value = 8048.7856
$ws.Range("I106").Value = 8048.7856
$ws.Range("K106").Value = 8048.7856
$ws.Range("M106").Value = -7417.7856

$ws.Range("H112").Value = 4899.5264
$ws.Range("J112").Value = 5096.222
$ws.Range("L112").Value = 15288.666
$ws.Range("N112").Value = -17504.666

$ws.Range("H132").Value = 1895.1464
$ws.Range("I132").Value = 1895.5658
$ws.Range("J132").Value = 1889.8334
$ws.Range("K132").Value = 5686.6974
$ws.Range("L132").Value = 5669.5002
$ws.Range("M132").Value = -3156.6974
$ws.Range("N132").Value = -10729.5002

$ws.Range("H137").Value = 2560.6667
$ws.Range("I137").Value = 2469.5
$ws.Range("J137").Value = 2674.625
$ws.Range("K137").Value = 7408.5
$ws.Range("L137").Value = 8023.875
$ws.Range("M137").Value = -4858.5
$ws.Range("N137").Value = -13123.875

$ws.Range("H138").Value = 3092.3142
$ws.Range("I138").Value = 1415.8667
$ws.Range("J138").Value = 4349.65
$ws.Range("K138").Value = 4247.6001
$ws.Range("L138").Value = 13048.95
$ws.Range("M138").Value = 892.3999000000003
$ws.Range("N138").Value = -23328.95

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6040.59
$ws.Range("I32").Value = 808.63635
$ws.Range("J32").Value = 34816.332
$ws.Range("K32").Value = 808.63635
$ws.Range("L32").Value = 34816.332
$ws.Range("M32").Value = -521.63635
$ws.Range("N32").Value = -35390.332

$ws.Range("H61").Value = 4890.9033
$ws.Range("I61").Value = 3337.238
$ws.Range("J61").Value = 8153.6
$ws.Range("K61").Value = 3337.238
$ws.Range("L61").Value = 8153.6
$ws.Range("M61").Value = -3125.238
$ws.Range("N61").Value = -8577.6

$ws.Range("H74").Value = 1599.775
$ws.Range("I74").Value = 1347.1389
$ws.Range("J74").Value = 3873.5
$ws.Range("K74").Value = 1347.1389
$ws.Range("L74").Value = 3873.5
$ws.Range("M74").Value = -473.1388999999999
$ws.Range("N74").Value = -5621.5

$ws.Range("H75").Value = 25173
$ws.Range("J75").Value = 25173
$ws.Range("L75").Value = 25173
$ws.Range("N75").Value = -26921

$ws.Range("H77").Value = 1599.775
$ws.Range("I77").Value = 1347.1389
$ws.Range("J77").Value = 3873.5
$ws.Range("K77").Value = 6735.6945
$ws.Range("L77").Value = 19367.5
$ws.Range("M77").Value = -2367.6945
$ws.Range("N77").Value = -28103.5

$ws.Range("H78").Value = 25173
$ws.Range("J78").Value = 25173
$ws.Range("L78").Value = 75519
$ws.Range("N78").Value = -84255

$ws.Range("H132").Value = 2939.9756
$ws.Range("I132").Value = 2464.121
$ws.Range("J132").Value = 4902.875
$ws.Range("K132").Value = 7392.363
$ws.Range("L132").Value = 14708.625
$ws.Range("M132").Value = -4862.363
$ws.Range("N132").Value = -19768.625

$ws.Range("H136").Value = 4890.9033
$ws.Range("I136").Value = 3337.238
$ws.Range("J136").Value = 8153.6
$ws.Range("K136").Value = 10011.714
$ws.Range("L136").Value = 24460.8
$ws.Range("M136").Value = -7461.714
$ws.Range("N136").Value = -29560.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 1088.4286
$ws.Range("I12").Value = 754.75
$ws.Range("J12").Value = 1533.3334
$ws.Range("K12").Value = 754.75
$ws.Range("L12").Value = 1533.3334
$ws.Range("M12").Value = -586.75
$ws.Range("N12").Value = -1869.3334

$ws.Range("H107").Value = 3152.5
$ws.Range("I107").Value = 2442.25
$ws.Range("J107").Value = 4099.5
$ws.Range("K107").Value = 2442.25
$ws.Range("L107").Value = 4099.5
$ws.Range("M107").Value = -522.25
$ws.Range("N107").Value = -7939.5

$ws.Range("H128").Value = 1999
$ws.Range("I128").Value = 1999
$ws.Range("K128").Value = 5997
$ws.Range("M128").Value = -3507

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H31").Value = 5388.881
$ws.Range("I31").Value = 4458.6206
$ws.Range("K31").Value = 4458.6206
$ws.Range("M31").Value = -4163.6206

$ws.Range("H34").Value = 5388.881
$ws.Range("I34").Value = 4458.6206
$ws.Range("K34").Value = 4458.6206
$ws.Range("M34").Value = -4256.6206

$ws.Range("H132").Value = 1035.3247
$ws.Range("I132").Value = 950.4925500000001
$ws.Range("J132").Value = 1603.7
$ws.Range("K132").Value = 2851.47765
$ws.Range("L132").Value = 4811.1
$ws.Range("M132").Value = -321.4776500000003
$ws.Range("N132").Value = -9871.1

$ws.Range("H134").Value = 1338.6438
$ws.Range("I134").Value = 1280.884
$ws.Range("K134").Value = 3842.652
$ws.Range("M134").Value = -1307.652

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 549.4815
$ws.Range("I5").Value = 447.4737
$ws.Range("J5").Value = 791.75
$ws.Range("K5").Value = 1342.4211
$ws.Range("L5").Value = 2375.25
$ws.Range("M5").Value = -1230.4211
$ws.Range("N5").Value = -2599.25

$ws.Range("H11").Value = 4062.3333
$ws.Range("I11").Value = 594
$ws.Range("J11").Value = 10999
$ws.Range("K11").Value = 1782
$ws.Range("L11").Value = 32997
$ws.Range("M11").Value = -1642
$ws.Range("N11").Value = -33277

$ws.Range("H46").Value = 23963.045
$ws.Range("I46").Value = 1223.8823
$ws.Range("J46").Value = 101276.2
$ws.Range("K46").Value = 3671.6469
$ws.Range("L46").Value = 303828.6
$ws.Range("M46").Value = -3580.6469
$ws.Range("N46").Value = -304010.6

$ws.Range("H135").Value = 549.4815
$ws.Range("I135").Value = 447.4737
$ws.Range("J135").Value = 791.75
$ws.Range("K135").Value = 4027.2633
$ws.Range("L135").Value = 7125.75
$ws.Range("M135").Value = -1492.2633
$ws.Range("N135").Value = -12195.75

$ws.Range("H136").Value = 6467.6665
$ws.Range("J136").Value = 8886.75
$ws.Range("L136").Value = 26660.25
$ws.Range("N136").Value = -36860.25

$ws.Range("H139").Value = 4035.0386
$ws.Range("I139").Value = 3207.2632
$ws.Range("J139").Value = 6281.857
$ws.Range("K139").Value = 9621.7896
$ws.Range("L139").Value = 18845.571
$ws.Range("M139").Value = -4481.7896
$ws.Range("N139").Value = -29125.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 29928
$ws.Range("J15").Value = 29928
$ws.Range("L15").Value = 29928
$ws.Range("N15").Value = -30504

$ws.Range("H39").Value = 25000
$ws.Range("J39").Value = 25000
$ws.Range("L39").Value = 25000
$ws.Range("N39").Value = -26064

$ws.Range("H81").Value = 29928
$ws.Range("J81").Value = 29928
$ws.Range("L81").Value = 29928
$ws.Range("N81").Value = -31924

$ws.Range("H84").Value = 29928
$ws.Range("J84").Value = 29928
$ws.Range("L84").Value = 89784
$ws.Range("N84").Value = -99768

$ws.Range("H106").Value = 25000
$ws.Range("J106").Value = 25000
$ws.Range("L106").Value = 25000
$ws.Range("N106").Value = -27524

$ws.Range("H126").Value = 6001.8
$ws.Range("J126").Value = 10000
$ws.Range("L126").Value = 30000
$ws.Range("N126").Value = -34940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()

$ws.Range("H76").Value = 22899
$ws.Range("I76").Value = 5898
$ws.Range("J76").Value = 39900
$ws.Range("K76").Value = 5898
$ws.Range("L76").Value = 39900
$ws.Range("M76").Value = -5560
$ws.Range("N76").Value = -40576

$ws.Range("H79").Value = 22899
$ws.Range("I79").Value = 5898
$ws.Range("J79").Value = 39900
$ws.Range("K79").Value = 5898
$ws.Range("L79").Value = 39900
$ws.Range("M79").Value = -4728
$ws.Range("N79").Value = -42240

$ws.Range("H93").Value = 798.4737
$ws.Range("I93").Value = 744.17645
$ws.Range("J93").Value = 1260
$ws.Range("K93").Value = 744.17645
$ws.Range("L93").Value = 1260
$ws.Range("M93").Value = 503.82355
$ws.Range("N93").Value = -3756

$ws.Range("H97").Value = 32000
$ws.Range("J97").Value = 32000
$ws.Range("L97").Value = 32000
$ws.Range("N97").Value = -33982

$ws.Range("H105").Value = 59966.668
$ws.Range("J105").Value = 59966.668
$ws.Range("L105").Value = 59966.668
$ws.Range("N105").Value = -66954.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 35275.5
$ws.Range("I64").Value = 26103
$ws.Range("J64").Value = 38333
$ws.Range("K64").Value = 26103
$ws.Range("L64").Value = 38333
$ws.Range("M64").Value = -25855
$ws.Range("N64").Value = -38829

$ws.Range("H67").Value = 35275.5
$ws.Range("I67").Value = 26103
$ws.Range("J67").Value = 38333
$ws.Range("K67").Value = 26103
$ws.Range("L67").Value = 38333
$ws.Range("M67").Value = -25245
$ws.Range("N67").Value = -40049

$ws.Range("H136").Value = 5654
$ws.Range("I136").Value = 6003.1177
$ws.Range("J136").Value = 4664.8335
$ws.Range("K136").Value = 18009.3531
$ws.Range("L136").Value = 13994.5005
$ws.Range("M136").Value = -15459.3531
$ws.Range("N136").Value = -19094.5005
